$wb = $excel.ActiveWorkbook

$sheetNames = @("JFK", "Regular", "Others")

$cValues = @{
    "JFK"     = @(40937, 49903, 77429)
    "Regular" = @(1582901, 1908848, 2320784)
    "Others"  = @(12792, 15274, 20292)
}
$dValues = @{
    "JFK"     = @(749144.01, 909913.16, 1409270.6)
    "Regular" = @(4076706.02, 4950829.53, 6326754.48)
    "Others"  = @(173169.86, 202062.68, 285632.9)
}
$eValues = @{
    "JFK"     = @(63641, 77372, 120217)
    "Regular" = @(2197129, 2658010, 3228655)
    "Others"  = @(16532, 19335, 26542)
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A1").Value = "year_month"

    $rows = @(2, 4, 6)
    for ($i = 0; $i -lt $rows.Length; $i++) {
        $r = $rows[$i]
        $ws.Cells.Item($r, 3).Value = $cValues[$name][$i]
        $ws.Cells.Item($r, 4).Value = $dValues[$name][$i]
        $ws.Cells.Item($r, 5).Value = $eValues[$name][$i]
    }
}
